$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 261.54544
$ws.Range("I33").Value = 256.42307
$ws.Range("J33").Value = 280.57144
$ws.Range("K33").Value = 256.42307
$ws.Range("L33").Value = 280.57144
$ws.Range("M33").Value = -27.42307
$ws.Range("N33").Value = -738.5714399999999
$ws.Range("H116").Value = 6109.091
$ws.Range("I116").Value = 1680
$ws.Range("J116").Value = 9800
$ws.Range("K116").Value = 1680
$ws.Range("L116").Value = 9800
$ws.Range("M116").Value = 1762
$ws.Range("N116").Value = -16684
$ws.Range("H132").Value = 1887075.2
$ws.Range("I132").Value = 2558.28
$ws.Range("K132").Value = 7674.84
$ws.Range("M132").Value = -5144.84
$ws.Range("H135").Value = 41609.68
$ws.Range("I135").Value = 53895.316
$ws.Range("J135").Value = 2705.1667
$ws.Range("K135").Value = 485057.844
$ws.Range("L135").Value = 24346.5003
$ws.Range("M135").Value = -482522.844
$ws.Range("N135").Value = -29416.5003
$ws.Range("H136").Value = 49869
$ws.Range("J136").Value = 49869
$ws.Range("L136").Value = 49869
$ws.Range("N136").Value = -60069
$ws.Range("H137").Value = 2634181.8
$ws.Range("I137").Value = 4349839.5
$ws.Range("J137").Value = 3506.6
$ws.Range("K137").Value = 13049518.5
$ws.Range("L137").Value = 10519.8
$ws.Range("M137").Value = -13046968.5
$ws.Range("N137").Value = -15619.8
$ws.Range("H138").Value = 1917773.8
$ws.Range("I138").Value = 1218.25
$ws.Range("J138").Value = 2489879.8
$ws.Range("K138").Value = 3654.75
$ws.Range("L138").Value = 7469639.399999999
$ws.Range("M138").Value = 1485.25
$ws.Range("N138").Value = -7479919.399999999
$ws.Range("H141").Value = 1337.2667
$ws.Range("I141").Value = 1289.9286
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 3869.7858
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 1310.2142
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 125251150
$ws.Range("I61").Value = 200201090
$ws.Range("J61").Value = 334600
$ws.Range("K61").Value = 200201090
$ws.Range("L61").Value = 334600
$ws.Range("M61").Value = -200200878
$ws.Range("N61").Value = -335024
$ws.Range("H74").Value = 14001151
$ws.Range("I74").Value = 17929562
$ws.Range("J74").Value = 251712.5
$ws.Range("K74").Value = 17929562
$ws.Range("L74").Value = 251712.5
$ws.Range("M74").Value = -17928688
$ws.Range("N74").Value = -253460.5
$ws.Range("H77").Value = 14001151
$ws.Range("I77").Value = 17929562
$ws.Range("J77").Value = 251712.5
$ws.Range("K77").Value = 89647810
$ws.Range("L77").Value = 1258562.5
$ws.Range("M77").Value = -89643442
$ws.Range("N77").Value = -1267298.5
$ws.Range("H88").Value = 5273.1304
$ws.Range("I88").Value = 2686.5557
$ws.Range("J88").Value = 6935.9287
$ws.Range("K88").Value = 2686.5557
$ws.Range("L88").Value = 6935.9287
$ws.Range("M88").Value = -2280.5557
$ws.Range("N88").Value = -7747.9287
$ws.Range("H91").Value = 5273.1304
$ws.Range("I91").Value = 2686.5557
$ws.Range("J91").Value = 6935.9287
$ws.Range("K91").Value = 2686.5557
$ws.Range("L91").Value = 6935.9287
$ws.Range("M91").Value = -1282.5557
$ws.Range("N91").Value = -9743.9287
$ws.Range("H132").Value = 108102.69
$ws.Range("I132").Value = 79841.30499999999
$ws.Range("J132").Value = 169335.67
$ws.Range("K132").Value = 239523.915
$ws.Range("L132").Value = 508007.01
$ws.Range("M132").Value = -236993.915
$ws.Range("N132").Value = -513067.01
$ws.Range("H136").Value = 125251150
$ws.Range("I136").Value = 200201090
$ws.Range("J136").Value = 334600
$ws.Range("K136").Value = 600603270
$ws.Range("L136").Value = 1003800
$ws.Range("M136").Value = -600600720
$ws.Range("N136").Value = -1008900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 594.75
$ws.Range("I22").Value = 419.85715
$ws.Range("J22").Value = 839.6
$ws.Range("K22").Value = 419.85715
$ws.Range("L22").Value = 839.6
$ws.Range("M22").Value = -246.85715
$ws.Range("N22").Value = -1185.6
$ws.Range("H134").Value = 4200.1797
$ws.Range("I134").Value = 3643.1
$ws.Range("J134").Value = 6057.1113
$ws.Range("K134").Value = 10929.3
$ws.Range("L134").Value = 18171.3339
$ws.Range("M134").Value = -8394.299999999999
$ws.Range("N134").Value = -23241.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 992.5
$ws.Range("I12").Value = 992.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 992.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -822.5
$ws.Range("N12").ClearContents()
$ws.Range("H31").Value = 2716.73
$ws.Range("I31").Value = 1274.96
$ws.Range("J31").Value = 3197.32
$ws.Range("K31").Value = 1274.96
$ws.Range("L31").Value = 3197.32
$ws.Range("M31").Value = -979.96
$ws.Range("N31").Value = -3787.32
$ws.Range("H34").Value = 2716.73
$ws.Range("I34").Value = 1274.96
$ws.Range("J34").Value = 3197.32
$ws.Range("K34").Value = 1274.96
$ws.Range("L34").Value = 3197.32
$ws.Range("M34").Value = -1072.96
$ws.Range("N34").Value = -3601.32
$ws.Range("H58").Value = 27030124
$ws.Range("I58").Value = 33336738
$ws.Range("J58").Value = 1770.8572
$ws.Range("K58").Value = 33336738
$ws.Range("L58").Value = 1770.8572
$ws.Range("M58").Value = -33336535
$ws.Range("N58").Value = -2176.8572
$ws.Range("H132").Value = 115332
$ws.Range("I132").Value = 3200
$ws.Range("J132").Value = 129348.5
$ws.Range("K132").Value = 9600
$ws.Range("L132").Value = 388045.5
$ws.Range("M132").Value = -7070
$ws.Range("N132").Value = -393105.5
$ws.Range("H134").Value = 136461.25
$ws.Range("I134").Value = 1922.5
$ws.Range("J134").Value = 271000
$ws.Range("K134").Value = 5767.5
$ws.Range("L134").Value = 813000
$ws.Range("M134").Value = -3232.5
$ws.Range("N134").Value = -818070
$ws.Range("H136").Value = 27030124
$ws.Range("I136").Value = 33336738
$ws.Range("J136").Value = 1770.8572
$ws.Range("K136").Value = 100010214
$ws.Range("L136").Value = 5312.571599999999
$ws.Range("M136").Value = -100007664
$ws.Range("N136").Value = -10412.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 56363.277
$ws.Range("I5").Value = 111693.78
$ws.Range("J5").Value = 1032.7778
$ws.Range("K5").Value = 335081.34
$ws.Range("L5").Value = 3098.3334
$ws.Range("M5").Value = -334969.34
$ws.Range("N5").Value = -3322.3334
$ws.Range("H122").Value = 1090.5862
$ws.Range("I122").Value = 507.25
$ws.Range("J122").Value = 1312.8096
$ws.Range("K122").Value = 4565.25
$ws.Range("L122").Value = 11815.2864
$ws.Range("M122").Value = -2115.25
$ws.Range("N122").Value = -16715.2864
$ws.Range("H132").Value = 1903.8462
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560
$ws.Range("H135").Value = 56363.277
$ws.Range("I135").Value = 111693.78
$ws.Range("J135").Value = 1032.7778
$ws.Range("K135").Value = 1005244.02
$ws.Range("L135").Value = 9295.0002
$ws.Range("M135").Value = -1002709.02
$ws.Range("N135").Value = -14365.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1132.8889
$ws.Range("I102").Value = 1132.6666
$ws.Range("J102").Value = 1133.3334
$ws.Range("K102").Value = 1132.6666
$ws.Range("L102").Value = 1133.3334
$ws.Range("M102").Value = 489.3334
$ws.Range("N102").Value = -4377.3334
$ws.Range("H122").Value = 4572.4
$ws.Range("I122").Value = 4458.7
$ws.Range("K122").Value = 13376.1
$ws.Range("M122").Value = -10926.1
$ws.Range("H132").Value = 90649.56
$ws.Range("I132").Value = 80523.53999999999
$ws.Range("J132").Value = 103813.4
$ws.Range("K132").Value = 241570.62
$ws.Range("L132").Value = 311440.2
$ws.Range("M132").Value = -239040.62
$ws.Range("N132").Value = -316500.2
$ws.Range("H133").Value = 52000
$ws.Range("J133").Value = 52000
$ws.Range("L133").Value = 52000
$ws.Range("N133").Value = -62120
$ws.Range("H135").Value = 52000
$ws.Range("J135").Value = 52000
$ws.Range("L135").Value = 52000
$ws.Range("N135").Value = -62140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2186.6
$ws.Range("I7").Value = 2128.5715
$ws.Range("J7").Value = 2999
$ws.Range("K7").Value = 2128.5715
$ws.Range("L7").Value = 2999
$ws.Range("M7").Value = -2016.5715
$ws.Range("N7").Value = -3223
$ws.Range("H126").Value = 2186.6
$ws.Range("I126").Value = 2128.5715
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 6385.7145
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -3915.7145
$ws.Range("N126").Value = -13937
$ws.Range("H132").Value = 59391.668
$ws.Range("I132").Value = 3322.4443
$ws.Range("J132").Value = 115460.89
$ws.Range("K132").Value = 9967.332900000001
$ws.Range("L132").Value = 346382.67
$ws.Range("M132").Value = -7437.332900000001
$ws.Range("N132").Value = -351442.67
$ws.Range("H136").Value = 149140
$ws.Range("I136").Value = 114976.664
$ws.Range("J136").Value = 177091.81
$ws.Range("K136").Value = 344929.992
$ws.Range("L136").Value = 531275.4299999999
$ws.Range("M136").Value = -342379.992
$ws.Range("N136").Value = -536375.4299999999
$ws.Range("H138").Value = 43000
$ws.Range("J138").Value = 43000
$ws.Range("L138").Value = 43000
$ws.Range("N138").Value = -53280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 336649.66
$ws.Range("I132").Value = 502450
$ws.Range("J132").Value = 253749.5
$ws.Range("K132").Value = 1507350
$ws.Range("L132").Value = 761248.5
$ws.Range("M132").Value = -1504820
$ws.Range("N132").Value = -766308.5
$ws.Range("H136").Value = 71249.414
$ws.Range("I136").Value = 48028.047
$ws.Range("J136").Value = 144230.86
$ws.Range("K136").Value = 144084.141
$ws.Range("L136").Value = 432692.58
$ws.Range("M136").Value = -141534.141
$ws.Range("N136").Value = -437792.58
